{"js": "// Replace the date line and each two-digit multiplication problem/answer\n// with its updated value, matching the OOXML diff exactly.\nconst replacements = [\n  [\"2026-01-22 Thursday\", \"2026-01-23 Friday\"],\n  [\"57\u00d785=4845\", \"62\u00d795=5890\"],\n  [\"87\u00d764=5568\", \"26\u00d718=468\"],\n  [\"59\u00d787=5133\", \"36\u00d754=1944\"],\n  [\"47\u00d736=1692\", \"93\u00d799=9207\"],\n  [\"14\u00d788=1232\", \"54\u00d719=1026\"],\n  [\"75\u00d799=7425\", \"31\u00d746=1426\"],\n  [\"52\u00d750=2600\", \"83\u00d779=6557\"],\n  [\"89\u00d786=7654\", \"77\u00d745=3465\"],\n  [\"49\u00d765=3185\", \"15\u00d717=255\"],\n  [\"41\u00d785=3485\", \"59\u00d768=4012\"],\n  [\"65\u00d760=3900\", \"23\u00d740=920\"],\n  [\"77\u00d796=7392\", \"98\u00d794=9212\"],\n  [\"63\u00d777=4851\", \"91\u00d773=6643\"],\n  [\"24\u00d751=1224\", \"50\u00d717=850\"],\n  [\"99\u00d795=9405\", \"92\u00d739=3588\"],\n  [\"36\u00d750=1800\", \"65\u00d731=2015\"],\n  [\"40\u00d724=960\", \"23\u00d780=1840\"],\n  [\"45\u00d725=1125\", \"26\u00d782=2132\"],\n  [\"73\u00d784=6132\", \"89\u00d785=7565\"],\n  [\"18\u00d734=612\", \"52\u00d734=1768\"],\n  [\"95\u00d766=6270\", \"44\u00d788=3872\"],\n  [\"62\u00d756=3472\", \"68\u00d779=5372\"],\n  [\"70\u00d765=4550\", \"71\u00d738=2698\"],\n  [\"84\u00d732=2688\", \"62\u00d769=4278\"],\n  [\"97\u00d759=5723\", \"98\u00d786=8428\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @{Old = \"2026-01-22 Thursday\"; New = \"2026-01-23 Friday\"},\n    @{Old = \"57\u00d785=4845\"; New = \"62\u00d795=5890\"},\n    @{Old = \"87\u00d764=5568\"; New = \"26\u00d718=468\"},\n    @{Old = \"59\u00d787=5133\"; New = \"36\u00d754=1944\"},\n    @{Old = \"47\u00d736=1692\"; New = \"93\u00d799=9207\"},\n    @{Old = \"14\u00d788=1232\"; New = \"54\u00d719=1026\"},\n    @{Old = \"75\u00d799=7425\"; New = \"31\u00d746=1426\"},\n    @{Old = \"52\u00d750=2600\"; New = \"83\u00d779=6557\"},\n    @{Old = \"89\u00d786=7654\"; New = \"77\u00d745=3465\"},\n    @{Old = \"49\u00d765=3185\"; New = \"15\u00d717=255\"},\n    @{Old = \"41\u00d785=3485\"; New = \"59\u00d768=4012\"},\n    @{Old = \"65\u00d760=3900\"; New = \"23\u00d740=920\"},\n    @{Old = \"77\u00d796=7392\"; New = \"98\u00d794=9212\"},\n    @{Old = \"63\u00d777=4851\"; New = \"91\u00d773=6643\"},\n    @{Old = \"24\u00d751=1224\"; New = \"50\u00d717=850\"},\n    @{Old = \"99\u00d795=9405\"; New = \"92\u00d739=3588\"},\n    @{Old = \"36\u00d750=1800\"; New = \"65\u00d731=2015\"},\n    @{Old = \"40\u00d724=960\"; New = \"23\u00d780=1840\"},\n    @{Old = \"45\u00d725=1125\"; New = \"26\u00d782=2132\"},\n    @{Old = \"73\u00d784=6132\"; New = \"89\u00d785=7565\"},\n    @{Old = \"18\u00d734=612\"; New = \"52\u00d734=1768\"},\n    @{Old = \"95\u00d766=6270\"; New = \"44\u00d788=3872\"},\n    @{Old = \"62\u00d756=3472\"; New = \"68\u00d779=5372\"},\n    @{Old = \"70\u00d765=4550\"; New = \"71\u00d738=2698\"},\n    @{Old = \"84\u00d732=2688\"; New = \"62\u00d769=4278\"},\n    @{Old = \"97\u00d759=5723\"; New = \"98\u00d786=8428\"}\n)\n\nforeach ($r in $replacements) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $r.Old\n    $find.Replacement.Text = $r.New\n    $find.Forward = $true\n    $find.Wrap = 1  # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $false, $false, $false, $find.Forward, $find.Wrap, $false, $find.Replacement.Text, 2)  # wdReplaceAll\n}\n"}
